# Modified reference data replacement logic
#
# Applies the target edits to the "Data" and "Reference" sheets:
#  - Data sheet: TestCase1 table row values switched from raw placeholder
#    text (tcXrYcZ) to real templated strings / numbers, Run Mode for the
#    3rd record flipped to "Y", number formats tweaked for numeric /
#    date columns, and a new 6th field column populated.
#  - Data sheet: TestCase2 table 2nd record columns now hold
#    ${userFirstName} / ${userLastName} placeholders.
#  - Reference sheet: grew from 2 placeholder rows to 4, plus a computed
#    "full name" row built with CONCATENATE().

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Data sheet
# ---------------------------------------------------------------------
$data = $wb.Worksheets.Item("Data")

# Row 3 (TestCase1, record 1): tweak number format of the decimal column
# and bump the date by one day, switching its display format too.
$data.Range("C3").NumberFormat = "0.00"
$data.Range("D3").Value = 43203
$data.Range("D3").NumberFormat = "dd/mm/yyyy"

# Row 4 (TestCase1, record 2): was all placeholder text, now real data.
$data.Range("B4").Value = 22
$data.Range("C4").Value = 45.89
$data.Range("C4").NumberFormat = "0.00"
$data.Range("D4").Value = "welcome `${userName}!"
$data.Range("E4").Value = "`${userName} last logged in at `${myDate}"
$data.Range("F4").Value = "tc1r2c5"

# Row 5 (TestCase1, record 3): Run Mode flips to Y, placeholder text
# replaced with real data (F5 keeps its original text).
$data.Range("A5").Value = "Y"
$data.Range("B5").Value = 33
$data.Range("C5").Value = 78.224
$data.Range("C5").NumberFormat = "0.00"
$data.Range("D5").Value = "welcome `${userFullName}!"
$data.Range("E5").Value = "`${myDate}"
$data.Range("F5").Value = "tc1r3c5"

# Row 10 (TestCase2, record 2): placeholder text replaced.
$data.Range("B10").Value = "`${userFirstName}"
$data.Range("C10").Value = "`${userLastName}"

# Column widths widened to fit the new, longer cell content.
$data.Columns("B").ColumnWidth = 15.59
$data.Columns("C").ColumnWidth = 15.25
$data.Columns("D").ColumnWidth = 24.76
$data.Columns("E").ColumnWidth = 36.59
$data.Columns("F").ColumnWidth = 6.25
$data.Columns("G").ColumnWidth = 6.25

$data.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# Reference sheet
# ---------------------------------------------------------------------
$ref = $wb.Worksheets.Item("Reference")

# Existing rows keep their values; just refresh the date cell's format
# and left-align the "swapnil" value cell to match the new rows below.
$ref.Range("B2").HorizontalAlignment = -4131
$ref.Range("B3").NumberFormat = "dd\-mm\-yyyy"
$ref.Range("B3").HorizontalAlignment = -4131

# New rows 4-6: first/last name placeholders plus a computed full name.
$ref.Range("A4").Value = "`${userFirstName}"
$ref.Range("B4").Value = "swapnil"
$ref.Range("B4").HorizontalAlignment = -4131

$ref.Range("A5").Value = "`${userLastName}"
$ref.Range("B5").Value = "sonar"
$ref.Range("B5").HorizontalAlignment = -4131

$ref.Range("A6").Value = "`${userFullName}"
$ref.Range("B6").Formula = "=CONCATENATE(B4,"" "",B5)"
$ref.Range("B6").HorizontalAlignment = -4131

$ref.Columns("A").ColumnWidth = 15.59
$ref.Columns("B").ColumnWidth = 12.42
